# Scheduled runner update: refresh cached market-board price snapshots
# (currentAveragePrice*/LevePrice*/LeveProfit*) across the per-job Leve
# profit sheets. Values below are the new snapshot figures; every other
# cell in each row is left untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 99600
$ws.Range("J3").Value = 99600
$ws.Range("L3").Value = 99600
$ws.Range("N3").Value = -99828

$ws.Range("H33").Value = 52644052
$ws.Range("I33").Value = 76925760
$ws.Range("J33").Value = 33683
$ws.Range("K33").Value = 76925760
$ws.Range("L33").Value = 33683
$ws.Range("M33").Value = -76925531
$ws.Range("N33").Value = -34141

$ws.Range("H102").Value = 99600
$ws.Range("J102").Value = 99600
$ws.Range("L102").Value = 99600
$ws.Range("N102").Value = -106090

$ws.Range("H103").Value = 2778452.8
$ws.Range("I103").Value = 11111111
$ws.Range("J103").Value = 900
$ws.Range("K103").Value = 33333333
$ws.Range("L103").Value = 2700
$ws.Range("M103").Value = -33332747
$ws.Range("N103").Value = -3872

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5255.0146
$ws.Range("I32").Value = 2623.544
$ws.Range("J32").Value = 18890.818
$ws.Range("K32").Value = 2623.544
$ws.Range("L32").Value = 18890.818
$ws.Range("M32").Value = -2336.544
$ws.Range("N32").Value = -19464.818

$ws.Range("H61").Value = 1527.8572
$ws.Range("I61").Value = 1260.3334
$ws.Range("K61").Value = 1260.3334
$ws.Range("M61").Value = -1048.3334

$ws.Range("H74").Value = 55556684
$ws.Range("I74").Value = 71429520
$ws.Range("J74").Value = 1757
$ws.Range("K74").Value = 71429520
$ws.Range("L74").Value = 1757
$ws.Range("M74").Value = -71428646
$ws.Range("N74").Value = -3505

$ws.Range("H77").Value = 55556684
$ws.Range("I77").Value = 71429520
$ws.Range("J77").Value = 1757
$ws.Range("K77").Value = 357147600
$ws.Range("L77").Value = 8785
$ws.Range("M77").Value = -357143232
$ws.Range("N77").Value = -17521

$ws.Range("H105").Value = 20370
$ws.Range("J105").Value = 20370
$ws.Range("L105").Value = 20370
$ws.Range("N105").Value = -27358

$ws.Range("H122").Value = 1728.5476
$ws.Range("I122").Value = 1460.7333
$ws.Range("K122").Value = 4382.199900000001
$ws.Range("M122").Value = -1932.199900000001

$ws.Range("H132").Value = 1456.6945
$ws.Range("I132").Value = 993.8261
$ws.Range("J132").Value = 2275.6155
$ws.Range("K132").Value = 2981.4783
$ws.Range("L132").Value = 6826.8465
$ws.Range("M132").Value = -451.4782999999998
$ws.Range("N132").Value = -11886.8465

$ws.Range("H136").Value = 1527.8572
$ws.Range("I136").Value = 1260.3334
$ws.Range("K136").Value = 3781.0002
$ws.Range("M136").Value = -1231.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 742.375
$ws.Range("I80").Value = 46
$ws.Range("J80").Value = 841.8570999999999
$ws.Range("K80").Value = 46
$ws.Range("L80").Value = 841.8570999999999
$ws.Range("M80").Value = 952
$ws.Range("N80").Value = -2837.8571

$ws.Range("H83").Value = 742.375
$ws.Range("I83").Value = 46
$ws.Range("J83").Value = 841.8570999999999
$ws.Range("K83").Value = 230
$ws.Range("L83").Value = 4209.2855
$ws.Range("M83").Value = 4762
$ws.Range("N83").Value = -14193.2855

$ws.Range("H100").Value = 29800
$ws.Range("J100").Value = 29800
$ws.Range("L100").Value = 29800
$ws.Range("N100").Value = -31964

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 126.2
$ws.Range("I7").Value = 125.388885
$ws.Range("J7").Value = 128.28572
$ws.Range("K7").Value = 125.388885
$ws.Range("L7").Value = 128.28572
$ws.Range("M7").Value = -12.388885
$ws.Range("N7").Value = -354.28572

$ws.Range("H22").Value = 31250448
$ws.Range("I22").Value = 45454860
$ws.Range("J22").Value = 740.8
$ws.Range("K22").Value = 45454860
$ws.Range("L22").Value = 740.8
$ws.Range("M22").Value = -45454510
$ws.Range("N22").Value = -1440.8

$ws.Range("H31").Value = 18522020
$ws.Range("I31").Value = 31252488
$ws.Range("J31").Value = 4974.273
$ws.Range("K31").Value = 31252488
$ws.Range("L31").Value = 4974.273
$ws.Range("M31").Value = -31252193
$ws.Range("N31").Value = -5564.273

$ws.Range("H34").Value = 18522020
$ws.Range("I34").Value = 31252488
$ws.Range("J34").Value = 4974.273
$ws.Range("K34").Value = 31252488
$ws.Range("L34").Value = 4974.273
$ws.Range("M34").Value = -31252286
$ws.Range("N34").Value = -5378.273

$ws.Range("H106").Value = 29115.5
$ws.Range("J106").Value = 29115.5
$ws.Range("L106").Value = 29115.5
$ws.Range("N106").Value = -31639.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 1285560.2
$ws.Range("I109").Value = 1889.75
$ws.Range("J109").Value = 3339433.2
$ws.Range("K109").Value = 5669.25
$ws.Range("L109").Value = 10018299.6
$ws.Range("M109").Value = -4629.25
$ws.Range("N109").Value = -10020379.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 93214.17999999999
$ws.Range("I80").Value = 2521.4285
$ws.Range("J80").Value = 251926.5
$ws.Range("K80").Value = 2521.4285
$ws.Range("L80").Value = 251926.5
$ws.Range("M80").Value = -1523.4285
$ws.Range("N80").Value = -253922.5

$ws.Range("H83").Value = 93214.17999999999
$ws.Range("I83").Value = 2521.4285
$ws.Range("J83").Value = 251926.5
$ws.Range("K83").Value = 12607.1425
$ws.Range("L83").Value = 1259632.5
$ws.Range("M83").Value = -7615.1425
$ws.Range("N83").Value = -1269616.5

$ws.Range("H102").Value = 32112.182
$ws.Range("I102").Value = 41214.96
$ws.Range("J102").Value = 3666
$ws.Range("K102").Value = 41214.96
$ws.Range("L102").Value = 3666
$ws.Range("M102").Value = -39592.96
$ws.Range("N102").Value = -6910

$ws.Range("H113").Value = 8497.929
$ws.Range("I113").Value = 12252.444
$ws.Range("J113").Value = 1739.8
$ws.Range("K113").Value = 12252.444
$ws.Range("L113").Value = 1739.8
$ws.Range("M113").Value = -10082.444
$ws.Range("N113").Value = -6079.8

$ws.Range("H122").Value = 2181.6
$ws.Range("I122").Value = 2154.8823
$ws.Range("J122").Value = 2333
$ws.Range("K122").Value = 6464.646900000001
$ws.Range("L122").Value = 6999
$ws.Range("M122").Value = -4014.646900000001
$ws.Range("N122").Value = -11899

$ws.Range("H126").Value = 1525.1
$ws.Range("I126").Value = 1244.5454
$ws.Range("J126").Value = 1868
$ws.Range("K126").Value = 3733.6362
$ws.Range("L126").Value = 5604
$ws.Range("M126").Value = -1263.6362
$ws.Range("N126").Value = -10544

$ws.Range("H132").Value = 5938.6787
$ws.Range("I132").Value = 7525.4443
$ws.Range("J132").Value = 3082.5
$ws.Range("K132").Value = 22576.3329
$ws.Range("L132").Value = 9247.5
$ws.Range("M132").Value = -20046.3329
$ws.Range("N132").Value = -14307.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 9930.571
$ws.Range("I68").Value = 16904
$ws.Range("J68").Value = 2957.1428
$ws.Range("K68").Value = 16904
$ws.Range("L68").Value = 2957.1428
$ws.Range("M68").Value = -16155
$ws.Range("N68").Value = -4455.1428

$ws.Range("H71").Value = 9930.571
$ws.Range("I71").Value = 16904
$ws.Range("J71").Value = 2957.1428
$ws.Range("K71").Value = 84520
$ws.Range("L71").Value = 14785.714
$ws.Range("M71").Value = -80776
$ws.Range("N71").Value = -22273.714

$ws.Range("H82").Value = 2658.4167
$ws.Range("I82").Value = 2625.5
$ws.Range("J82").Value = 2674.875
$ws.Range("K82").Value = 2625.5
$ws.Range("L82").Value = 2674.875
$ws.Range("M82").Value = -2264.5
$ws.Range("N82").Value = -3396.875

$ws.Range("H85").Value = 2658.4167
$ws.Range("I85").Value = 2625.5
$ws.Range("J85").Value = 2674.875
$ws.Range("K85").Value = 2625.5
$ws.Range("L85").Value = 2674.875
$ws.Range("M85").Value = -1377.5
$ws.Range("N85").Value = -5170.875

$ws.Range("H93").Value = 1551.7812
$ws.Range("I93").Value = 1182.9615
$ws.Range("K93").Value = 1182.9615
$ws.Range("M93").Value = 65.03850000000011

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3000
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 3000
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 3000
$ws.Range("M96").ClearContents()  # NQ price now equals HQ price; LeveProfitNQ no longer applies
$ws.Range("N96").Value = -5746

$ws.Range("H122").Value = 1478.5625
$ws.Range("I122").Value = 1200.25
$ws.Range("J122").Value = 1756.875
$ws.Range("K122").Value = 3600.75
$ws.Range("L122").Value = 5270.625
$ws.Range("M122").Value = -1150.75
$ws.Range("N122").Value = -10170.625
